$d = $word.ActiveDocument

# --- Change 1: "UC6) Attack" -> "UC6) Attack enemy" (first occurrence only,
#     the plain use-case list entry, not the later bold/underlined heading) ---
$uc6 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "UC6) Attack`r") {
        $uc6 = $p
        break
    }
}
$uc6.Range.InsertAfter(" enemy")

# --- Change 2: add a new "UC9) Flee enemy" paragraph right after the
#     "UC8) Item-shop transaction" entry in that same list ---
$uc8 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "UC8) Item-shop transaction`r") {
        $uc8 = $p
        break
    }
}
$uc8.Range.InsertParagraphAfter()
$uc9 = $uc8.Next()
$uc9.Range.InsertAfter("UC9) Flee enemy")
